# ============================================================
# chore: update Sheets via scheduled runner
#
# Refreshes the live-market-derived columns (H:N) on the per-job
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR):
#   H currentAveragePrice      K LevePriceNQ
#   I currentAveragePriceNQ    L LevePriceHQ
#   J currentAveragePriceHQ    M LeveProfitNQ
#                              N LeveProfitHQ
#
# Only the rows whose underlying market data changed since the
# last run are touched; other rows/columns are left untouched.
# ============================================================

$wb = $excel.ActiveWorkbook

# Writes one refreshed leve-row. Any argument equal to the
# "__SKIP__" sentinel leaves that column exactly as-is; $null
# clears the column (the source no longer reports a value there).
function Set-LeveRow {
    param($ws, $Row, $H, $I, $J, $K, $L, $M, $N)
    $values = [ordered]@{ H = $H; I = $I; J = $J; K = $K; L = $L; M = $M; N = $N }
    foreach ($col in $values.Keys) {
        $v = $values[$col]
        if ($v -eq "__SKIP__") { continue }
        $cell = $ws.Range("$col$Row")
        if ($v -eq $null) {
            $cell.ClearContents()
        } else {
            $cell.Value = $v
        }
    }
}

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
Set-LeveRow $ws 18 5000 0 5000 0 5000 $null -5568
Set-LeveRow $ws 62 6340.3335 2871.4285 9375.625 2871.4285 9375.625 -2247.4285 -10623.625
Set-LeveRow $ws 64 7512.7617 "__SKIP__" 8252.532999999999 "__SKIP__" 8252.532999999999 "__SKIP__" -8748.532999999999
Set-LeveRow $ws 65 6340.3335 2871.4285 9375.625 14357.1425 46878.125 -11237.1425 -53118.125
Set-LeveRow $ws 67 7512.7617 "__SKIP__" 8252.532999999999 "__SKIP__" 8252.532999999999 "__SKIP__" -9968.532999999999
Set-LeveRow $ws 74 7282.0386 5859.421 "__SKIP__" 5859.421 "__SKIP__" -4923.421 "__SKIP__"
Set-LeveRow $ws 77 7282.0386 5859.421 "__SKIP__" 29297.105 "__SKIP__" -24617.105 "__SKIP__"
Set-LeveRow $ws 106 5974 4807.857 "__SKIP__" 4807.857 "__SKIP__" -4176.857 "__SKIP__"
Set-LeveRow $ws 116 7871.5 5491.3335 "__SKIP__" 5491.3335 "__SKIP__" -2049.3335 "__SKIP__"
Set-LeveRow $ws 125 2503.9333 1551.4445 "__SKIP__" 13963.0005 "__SKIP__" -11503.0005 "__SKIP__"
Set-LeveRow $ws 137 3915.7896 3399.6667 4380.3 10199.0001 13140.9 -7649.000100000001 -18240.9

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
Set-LeveRow $ws 32 5082.1064 3778.6155 11436.625 3778.6155 11436.625 -3491.6155 -12010.625
Set-LeveRow $ws 74 25646316 41670296 "__SKIP__" 41670296 "__SKIP__" -41669422 "__SKIP__"
Set-LeveRow $ws 77 25646316 41670296 "__SKIP__" 208351480 "__SKIP__" -208347112 "__SKIP__"
Set-LeveRow $ws 97 1045 878.7 "__SKIP__" 878.7 "__SKIP__" -382.7 "__SKIP__"
Set-LeveRow $ws 98 0 "__SKIP__" 0 "__SKIP__" 0 "__SKIP__" $null
Set-LeveRow $ws 102 2877 2806.8572 3122.5 2806.8572 3122.5 -1184.8572 -6366.5
Set-LeveRow $ws 109 98799.25 "__SKIP__" 98799.25 "__SKIP__" 98799.25 "__SKIP__" -101573.25
Set-LeveRow $ws 122 3598.1667 2609.7144 4227.1816 7829.1432 12681.5448 -5379.1432 -17581.5448

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
Set-LeveRow $ws 86 3762.3333 2572.5334 6736.8335 2572.5334 6736.8335 -1449.5334 -8982.833500000001
Set-LeveRow $ws 89 3762.3333 2572.5334 6736.8335 12862.667 33684.1675 -7246.666999999999 -44916.1675
Set-LeveRow $ws 94 1938.2941 1227.0769 4249.75 1227.0769 4249.75 -776.0769 -5151.75

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
Set-LeveRow $ws 6 837.75 733.3333 "__SKIP__" 733.3333 "__SKIP__" -620.3333 "__SKIP__"
Set-LeveRow $ws 17 820 820 "__SKIP__" 820 "__SKIP__" -646 "__SKIP__"

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
Set-LeveRow $ws 131 6946155.5 "__SKIP__" 5118854.5 "__SKIP__" 15356563.5 "__SKIP__" -15366643.5
Set-LeveRow $ws 138 2992.8572 1825 10000 5475 30000 -335 -40280
Set-LeveRow $ws 139 2965.7837 2199.158 "__SKIP__" 6597.474 "__SKIP__" -1457.474 "__SKIP__"

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
Set-LeveRow $ws 119 73684 "__SKIP__" 73684 "__SKIP__" 73684 "__SKIP__" -83360
Set-LeveRow $ws 121 0 "__SKIP__" 0 "__SKIP__" 0 "__SKIP__" $null
Set-LeveRow $ws 122 5128.4287 4225.5356 6934.2144 12676.6068 20802.6432 -10226.6068 -25702.6432
Set-LeveRow $ws 124 78643.336 "__SKIP__" 78643.336 "__SKIP__" 78643.336 "__SKIP__" -88463.336

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
Set-LeveRow $ws 16 7981.5454 3655.3333 "__SKIP__" 3655.3333 "__SKIP__" -3485.3333 "__SKIP__"
Set-LeveRow $ws 55 2046.9546 647.0909 "__SKIP__" 647.0909 "__SKIP__" -474.0909 "__SKIP__"
Set-LeveRow $ws 61 3486.8333 2617.7273 5876.875 2617.7273 5876.875 -2415.7273 -6280.875
Set-LeveRow $ws 82 14541 10198.5 17146.5 10198.5 17146.5 -9837.5 -17868.5
Set-LeveRow $ws 85 14541 10198.5 17146.5 10198.5 17146.5 -8950.5 -19642.5
Set-LeveRow $ws 93 1479.2 "__SKIP__" 2000 "__SKIP__" 2000 "__SKIP__" -4496
Set-LeveRow $ws 113 3486.8333 2617.7273 5876.875 2617.7273 5876.875 -447.7273 -10216.875
Set-LeveRow $ws 135 70464.336 "__SKIP__" 70464.336 "__SKIP__" 70464.336 "__SKIP__" -80604.336

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
Set-LeveRow $ws 14 0 0 "__SKIP__" 0 "__SKIP__" $null "__SKIP__"
Set-LeveRow $ws 81 2709.8 1346.1333 "__SKIP__" 2692.2666 "__SKIP__" -1631.2666 "__SKIP__"
Set-LeveRow $ws 84 2709.8 1346.1333 "__SKIP__" 13461.333 "__SKIP__" -8157.332999999999 "__SKIP__"
Set-LeveRow $ws 96 66 66 0 66 0 1307 $null
Set-LeveRow $ws 132 3775.1562 2504.72 "__SKIP__" 7514.16 "__SKIP__" -4984.16 "__SKIP__"

Write-Host "Scheduled runner: refreshed 227 cells across 8 sheets"
